# Update the worksheet date and regenerate the 25 division problems/answers.
# Cell-targeted assignment is used for the table entries (instead of a
# global Find/Replace) because several old values coincide with other
# cells' new values, which would corrupt a naive sequential replace-all.

$d = $word.ActiveDocument

# Header date line.
$d.Content.Find.Execute("2026-01-29 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2026-01-30 Friday", 2)

$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "23÷5=4, 3"
$t.Cell(1, 2).Range.Text = "21÷5=4, 1"
$t.Cell(1, 3).Range.Text = "27÷9=3, 0"
$t.Cell(1, 4).Range.Text = "37÷8=4, 5"
$t.Cell(1, 5).Range.Text = "30÷6=5, 0"

# Row 5
$t.Cell(5, 1).Range.Text = "64÷8=8, 0"
$t.Cell(5, 2).Range.Text = "22÷5=4, 2"
$t.Cell(5, 3).Range.Text = "44÷3=14, 2"
$t.Cell(5, 4).Range.Text = "51÷6=8, 3"
$t.Cell(5, 5).Range.Text = "89÷8=11, 1"

# Row 9
$t.Cell(9, 1).Range.Text = "75÷2=37, 1"
$t.Cell(9, 2).Range.Text = "47÷3=15, 2"
$t.Cell(9, 3).Range.Text = "72÷8=9, 0"
$t.Cell(9, 4).Range.Text = "42÷5=8, 2"
$t.Cell(9, 5).Range.Text = "70÷5=14, 0"

# Row 13
$t.Cell(13, 1).Range.Text = "19÷6=3, 1"
$t.Cell(13, 2).Range.Text = "97÷5=19, 2"
$t.Cell(13, 3).Range.Text = "61÷3=20, 1"
$t.Cell(13, 4).Range.Text = "70÷9=7, 7"
$t.Cell(13, 5).Range.Text = "83÷7=11, 6"

# Row 17
$t.Cell(17, 1).Range.Text = "70÷7=10, 0"
$t.Cell(17, 2).Range.Text = "99÷5=19, 4"
$t.Cell(17, 3).Range.Text = "61÷7=8, 5"
$t.Cell(17, 4).Range.Text = "66÷9=7, 3"
$t.Cell(17, 5).Range.Text = "99÷4=24, 3"

Write-Output "edit complete"
